$wb = $excel.ActiveWorkbook
$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# Add one "is this value missing" row per numeric variable, recording that
# "NA" is used as the missing-value marker.
$wsCategories.Range("A2").Value = "age_years"
$wsCategories.Range("A3").Value = "age_months"
$wsCategories.Range("A4").Value = "height_"
$wsCategories.Range("A5").Value = "height_age"
$wsCategories.Range("A6").Value = "weight_"
$wsCategories.Range("A7").Value = "weight_age"

$wsCategories.Range("B2:B7").Value = "NA"
$wsCategories.Range("D2:D7").Value = "NA"
$wsCategories.Range("C2:C7").Value = $true

# Match the row height used throughout the rest of the sheet.
$wsCategories.Rows.Item(2).RowHeight = 15
$wsCategories.Rows.Item(3).RowHeight = 15
$wsCategories.Rows.Item(4).RowHeight = 15
$wsCategories.Rows.Item(5).RowHeight = 15
$wsCategories.Rows.Item(6).RowHeight = 15
$wsCategories.Rows.Item(7).RowHeight = 15

# The "Categories" sheet header: variable | name | isMissing | label
# ("missing" is renamed to "isMissing")
$wsCategories.Range("C1").Value = "isMissing"

# Match the per-variable cell formatting already used on the "Variables"
# sheet for these same variables.
$wsVariables.Range("A5").Copy()
$wsCategories.Range("A3").PasteSpecial(-4122)
$wsVariables.Range("A6").Copy()
$wsCategories.Range("A4").PasteSpecial(-4122)
$wsVariables.Range("A7").Copy()
$wsCategories.Range("A5").PasteSpecial(-4122)
$wsVariables.Range("A8").Copy()
$wsCategories.Range("A6").PasteSpecial(-4122)
$wsVariables.Range("A9").Copy()
$wsCategories.Range("A7").PasteSpecial(-4122)

# Leave the selections where the edit did.
$wsVariables.Range("A5:A9").Select()
$wsCategories.Range("B3").Select()
